$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove all existing hyperlinks first. Inserting a row below shifts the
# cell text/values but NOT the hyperlink anchors, so they would otherwise
# end up pointing at the wrong cells. We rebuild every hyperlink from scratch
# afterwards once all the data is in its final place.
$ws.Hyperlinks.Delete()

# Insert a new row at the top of the data table (row 2). This shifts the
# existing rows 2-95 down to rows 3-96, carrying their values/styles along.
$ws.Rows.Item(2).Insert()

# Populate the newly inserted row 2 with the latest price entry.
# Column A/E hold dates stored as literal text (matching the rest of the
# sheet), so force a text number format before assigning them to avoid the
# values being auto-converted into date serials.
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "09-11-2025"
$ws.Range("B2").Value = "ALUMINIUM INGOT"
$ws.Range("C2").Value = "IE07"
$ws.Range("D2").Value = 297.15
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "01-11-2025"
$ws.Range("F2").Value = "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-01-11-2025.pdf"

# Copy the formatting from row 3 (a pre-existing data row) into row 2 so the
# new row matches the rest of the table (this also clears the "@" text
# format override above, restoring the normal General-format style).
$ws.Range("A3:F3").Copy()
$ws.Range("A2:F2").PasteSpecial(-4122)

# Re-create the hyperlink for every row in column F (the row insert above
# does not move hyperlink anchors, so they must all be rebuilt from scratch
# to point at the right row now that the table has shifted).
$ws.Hyperlinks.Add($ws.Range("F2"), "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-01-11-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F3"), "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-01-11-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F4"), "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-01-11-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F5"), "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-01-11-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F6"), "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-01-11-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F7"), "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-01-11-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F8"), "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-01-11-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F9"), "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-01-11-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F10"), "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-01-11-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F11"), "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-30-10-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F12"), "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-30-10-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F13"), "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-25-10-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F14"), "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-25-10-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F15"), "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-25-10-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F16"), "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-25-10-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F17"), "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-25-10-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F18"), "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-17-10-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F19"), "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-17-10-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F20"), "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-17-10-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F21"), "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-17-10-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F22"), "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-17-10-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F23"), "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-17-10-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F24"), "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-17-10-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F25"), "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-17-10-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F26"), "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-14-10-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F27"), "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-14-10-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F28"), "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-14-10-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F29"), "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-09-10-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F30"), "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-09-10-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F31"), "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-09-10-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F32"), "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-09-10-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F33"), "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-09-10-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F34"), "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-01-10-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F35"), "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-01-10-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F36"), "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-01-10-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F37"), "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-01-10-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F38"), "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-01-10-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F39"), "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-01-10-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F40"), "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-01-10-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F41"), "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-01-10-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F42"), "https://nalcoindia.com/wp-content/uploads/2025/09/INGOT-30-09-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F43"), "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-25-09-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F44"), "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-25-09-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F45"), "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-25-09-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F46"), "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-25-09-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F47"), "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-25-09-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F48"), "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-20-09-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F49"), "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-20-09-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F50"), "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-20-09-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F51"), "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-20-09-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F52"), "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-20-09-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F53"), "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-17-09-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F54"), "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-17-09-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F55"), "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-17-09-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F56"), "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-01-09-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F57"), "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-01-09-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F58"), "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-01-09-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F59"), "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-01-09-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F60"), "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-01-09-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F61"), "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-01-09-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F62"), "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-01-09-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F63"), "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-01-09-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F64"), "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-01-09-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F65"), "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-01-09-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F66"), "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-01-09-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F67"), "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-01-09-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F68"), "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-01-09-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F69"), "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-01-09-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F70"), "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-01-09-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F71"), "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-01-09-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F72"), "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-28-08-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F73"), "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-28-08-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F74"), "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-28-08-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F75"), "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-28-08-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F76"), "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-21-08-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F77"), "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-21-08-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F78"), "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-21-08-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F79"), "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-21-08-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F80"), "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-21-08-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F81"), "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-21-08-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F82"), "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-21-08-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F83"), "https://nalcoindia.com/wp-content/uploads/2025/08/INGOT-15-08-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F84"), "https://nalcoindia.com/wp-content/uploads/2025/08/INGOT-15-08-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F85"), "https://nalcoindia.com/wp-content/uploads/2025/08/INGOT-15-08-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F86"), "https://nalcoindia.com/wp-content/uploads/2025/08/INGOT-15-08-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F87"), "https://nalcoindia.com/wp-content/uploads/2025/08/INGOT-15-08-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F88"), "https://nalcoindia.com/wp-content/uploads/2025/08/INGOT-15-08-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F89"), "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-07-08-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F90"), "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-07-08-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F91"), "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-07-08-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F92"), "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-07-08-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F93"), "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-07-08-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F94"), "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-07-08-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F95"), "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-07-08-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F96"), "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-07-08-2025.pdf") | Out-Null
